# TDD clean-up of the RegisterUserData sheet: drop the fabricated sample
# records (Joseph Jackson / Joseph Dillon with email, password, interest,
# gender and hobby columns) and the test hyperlinks, leaving just two
# columns ("Name"/"State") with a single placeholder row, ready to be
# filled in at runtime by the Random-string test data generator.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RegisterUserData")
$ws.Activate()

# Drop the mailto: hyperlinks that were attached to the old Password column.
$ws.Hyperlinks.Delete()

# Remove the third data row entirely (shifts nothing else up, just gone).
$ws.Rows(3).Delete()

# Wipe all remaining cell content/formatting (columns C:G, old column
# widths, the old wrap/indent style) so the sheet goes back to plain
# defaults before we re-populate it.
$ws.Cells.Clear()
$ws.Columns("A:B").ClearFormats()
$ws.Rows("1:2").AutoFit()

# Re-seed just the two columns that remain.
$ws.Range("A1").Value = "Name"
$ws.Range("B1").Value = "State"
$ws.Range("A2").Value = "Darryal"
$ws.Range("B2").Value = "Goa"

$ws.Range("C2").Select() | Out-Null
